$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E (rows 2-51) to Text format so that numeric-looking
# strings (e.g. "1.003", "24.481.76") are stored as text, not auto-converted
# to numbers. We restore the style afterwards so no visible style/number
# format is left applied to the cells (matches original un-styled cells).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '24.481.76'
$ws.Range('E2').Value = '  -0.32%  '
$ws.Range('D3').Value = '1.658.70'
$ws.Range('E3').Value = '  -2.45%  '
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').Value = '307.44'
$ws.Range('E5').Value = '  -0.12%  '
$ws.Range('D6').Value = '0.9980'
$ws.Range('E6').Value = '  -0.19%  '
$ws.Range('D7').Value = '0.3623'
$ws.Range('E7').Value = '  -2.87%  '
$ws.Range('D8').Value = '47.40'
$ws.Range('E8').Value = '  -2.67%  '
$ws.Range('D9').Value = '0.3259'
$ws.Range('E9').Value = '  -4.73%  '
$ws.Range('E10').Value = '  -4.45%  '
$ws.Range('D11').Value = '0.06960'
$ws.Range('E11').Value = '  -6.24%  '
$ws.Range('D12').Value = '0.9997'
$ws.Range('E12').Value = '  -0.05%  '
$ws.Range('D13').Value = '5.915'
$ws.Range('E13').Value = '  -4.60%  '
$ws.Range('D14').Value = '19.34'
$ws.Range('E14').Value = '  -6.93%  '
$ws.Range('D15').Value = '6.594'
$ws.Range('E15').Value = '  -4.23%  '
$ws.Range('D16').Value = '1.655.45'
$ws.Range('E16').Value = '  -2.72%  '
$ws.Range('D17').Value = '0.00001044'
$ws.Range('E17').Value = '  -6.44%  '
$ws.Range('D18').Value = '0.06527'
$ws.Range('E18').Value = '  -2.34%  '
$ws.Range('D19').Value = '0.9981'
$ws.Range('E19').Value = '  -0.25%  '
$ws.Range('D20').Value = '76.49'
$ws.Range('E20').Value = '  -7.83%  '
$ws.Range('D21').Value = '5.907'
$ws.Range('E21').Value = '  -6.57%  '
$ws.Range('D22').Value = '15.71'
$ws.Range('D23').Value = '12.58'
$ws.Range('E23').Value = '  -3.75%  '
$ws.Range('D24').Value = '24.474.59'
$ws.Range('E24').Value = '  -0.29%  '
$ws.Range('D25').Value = '2.456'
$ws.Range('E25').Value = '  +1.57%  '
$ws.Range('E26').Value = '  -16.60%  '
$ws.Range('D27').Value = '146.76'
$ws.Range('E27').Value = '  -1.55%  '
$ws.Range('D28').Value = '18.45'
$ws.Range('D29').Value = '1.842.21'
$ws.Range('E29').Value = '  -2.43%  '
$ws.Range('D30').Value = '124.39'
$ws.Range('E30').Value = '  -4.86%  '
$ws.Range('D31').Value = '1.190'
$ws.Range('E31').Value = '  +1.46%  '
$ws.Range('E32').Value = '  -3.63%  '
$ws.Range('D33').Value = '5.598'
$ws.Range('E33').Value = '  -16.15%  '
$ws.Range('D34').Value = '0.08353'
$ws.Range('E34').Value = '  -4.41%  '
$ws.Range('E35').Value = '  -4.29%  '
$ws.Range('D36').Value = '12.36'
$ws.Range('E36').Value = '  -8.26%  '
$ws.Range('D37').Value = '5.198'
$ws.Range('E37').Value = '  -5.16%  '
$ws.Range('D38').Value = '0.06057'
$ws.Range('E38').Value = '  -6.69%  '
$ws.Range('B39').Value = 'Algorand'
$ws.Range('C39').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D39').Value = '0.2055'
$ws.Range('E39').Value = '  -6.81%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = '0.02189'
$ws.Range('E40').Value = '  -7.32%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').Value = '1.202'
$ws.Range('E41').Value = '  -5.45%  '
$ws.Range('D42').Value = '8.150'
$ws.Range('E42').Value = '  -8.34%  '
$ws.Range('D43').Value = '0.9987'
$ws.Range('E43').Value = '  -0.04%  '
$ws.Range('E44').Value = '  -7.61%  '
$ws.Range('D45').Value = '3.734'
$ws.Range('E45').Value = '  -1.59%  '
$ws.Range('D46').Value = '12.63'
$ws.Range('E46').Value = '  -8.17%  '
$ws.Range('D47').Value = '0.5600'
$ws.Range('E47').Value = '  -7.41%  '
$ws.Range('E48').Value = '  -4.98%  '
$ws.Range('D49').Value = '1.936'
$ws.Range('E49').Value = '  -7.78%  '
$ws.Range('D50').Value = '0.06920'
$ws.Range('E50').Value = '  -4.32%  '
$ws.Range('D51').Value = '74.02'
$ws.Range('E51').Value = '  -5.95%  '

# Restore original (unstyled) appearance for the Text-formatted range.
$ws.Range("D2:E51").Style = "Normal"
